# PROGRESO.xlsx edit script
# "se agregaron todos los instrumentos existentes"
# Adds a new "N of inst" column (G) to the instrument table, fills in
# previously-empty ADDED/WAV/PREVIEW cells for several instruments, and
# restructures the Organ (rows 10-12) block to match the Drum (rows 7-9) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in previously blank ADDED / WAV / PREVIEW values for several
#    instruments (Guitarra Electrica, Drum, Banjo).
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = 1

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("C5").Value = 1

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# ---------------------------------------------------------------------------
# 2. Organ block (rows 10-12) gets the same merged layout as the Drum block
#    (rows 7-9): D/E/F columns merged across the 3 rows and F10's formula
#    expanded to average across all three ADDED values.
# ---------------------------------------------------------------------------
$ws.Range("D10:D12").Merge()
$ws.Range("E10:E12").Merge()

$ws.Range("F10").Formula = "=C10/9 + C11/9 + C12/9 + D10/3 + E10/3"
$ws.Range("F10:F12").Merge()

# ---------------------------------------------------------------------------
# 3. New column G "N of inst" -- numbers every instrument in the order the
#    instruments were added to the project, with a running count total.
# ---------------------------------------------------------------------------
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1").Value = "N" + [char]176 + " of inst"

$ws.Range("D2").Copy()
$ws.Range("G2:G17").PasteSpecial(-4122) | Out-Null

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 9
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = ""
$ws.Range("G6").Value = 8
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("G13").Value = 3
$ws.Range("G14").Value = 4
$ws.Range("G15").Value = 5
$ws.Range("G16").Value = 6

$ws.Range("G4:G5").Merge()
$ws.Range("G7:G9").Merge()
$ws.Range("G10:G12").Merge()

$ws.Range("G17").Formula = "=MAX(G2:G16)+1"
$ws.Range("F17").Copy()
$ws.Range("G17").PasteSpecial(-4124) | Out-Null
$ws.Range("G17").Interior.Color = $ws.Range("D3").Interior.Color
$ws.Range("G17").Formula = "=MAX(G2:G16)+1"

# ---------------------------------------------------------------------------
# 4. Cosmetic: move the active selection like the author left it.
# ---------------------------------------------------------------------------
$ws.Range("I17").Select() | Out-Null
